$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.102621674537659
$ws.Range("B1").Value = 1.99652087688446
$ws.Range("C1").Value = 9.313411712646484
$ws.Range("D1").Value = 2.404149293899536
$ws.Range("E1").Value = 1.29102349281311
